$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Four new glucose/Hexose fragment rows (140-143), following the
# --- existing pattern of rows in this fragment-rules sheet ---

# name (col A)
$ws.Range("A140").Value = "[M+H-Hexose-H2O-CH4]+"
$ws.Range("A141").Value = "[M+H-Hexose-H2O-CH3OH]+"
$ws.Range("A142").Value = "[M+H-Hexose-H2O-C2H4O]+"
$ws.Range("A143").Value = "[M+H-Hexose-H2O-C2H4O2]+"

# nmol (col B)
$ws.Range("B140:B143").Value = 1

# charge (col C)
$ws.Range("C140:C143").Value = 1

# massdiff (col D) - formulas referencing the glucose parent row (139)
# and the corresponding neutral-loss rows
$ws.Range("D140").Formula = "=D139+D79+1.0073"
$ws.Range("D141").Formula = "=D139+D82+1.0073"
$ws.Range("D142").Formula = "=D139+D69+1.0073"
$ws.Range("D143").Formula = "=D139+D98+1.0073"

# oidscore (col E)
$ws.Range("E140").Value = 174
$ws.Range("E141").Value = 175
$ws.Range("E142").Value = 176
$ws.Range("E143").Value = 177

# quasi (col F)
$ws.Range("F140:F143").Value = 0

# ips (col G)
$ws.Range("G140:G143").Value = 0.5

# Match the fill/formatting used by the rest of this fragment block
# (row 139 is the adjacent, identically-styled row).
$ws.Range("A139:G139").Copy()
$ws.Range("A140:G143").PasteSpecial(-4122)

# Restore the selection state recorded in the workbook after the edit
$ws.Range("B138").Select()
